# Adding proportion calcs to bootstrapping
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E (proportion drinking) - numeric rows
$ws.Range("E2").Value = 0.105657
$ws.Range("E4").Value = 0.09945900000000001
$ws.Range("E6").Value = 0.068439
$ws.Range("E8").Value = 0.061954
$ws.Range("E10").Value = 0.054771
$ws.Range("E12").Value = 0.057475
$ws.Range("E14").Value = 0.050718

# Column C (theta SE) - text rows, unchanged cell refs but string pool values differ
$ws.Range("C3").Value = "(0.5)"
$ws.Range("C5").Value = "(0.43)"
$ws.Range("C7").Value = "(0.13)"
$ws.Range("C9").Value = "(0.1)"
$ws.Range("C11").Value = "(2.02)"
$ws.Range("C13").Value = "(1.09)"
$ws.Range("C15").Value = "(0.22)"

# Column D (lambda SE) - text rows
$ws.Range("D3").Value = "(0.12)"
$ws.Range("D5").Value = "(0.15)"
$ws.Range("D7").Value = "(0.17)"
$ws.Range("D9").Value = "(0.1)"
$ws.Range("D11").Value = "(2.16)"
$ws.Range("D13").Value = "(0.66)"
$ws.Range("D15").Value = "(0.31)"

# Column E (proportion drinking SE) - text rows
$ws.Range("E3").Value = "(0.00002)"
$ws.Range("E5").Value = "(0.00002)"
$ws.Range("E7").Value = "(0.00000)"
$ws.Range("E9").Value = "(0.00000)"
$ws.Range("E11").Value = "(0.00002)"
$ws.Range("E13").Value = "(0.00001)"
$ws.Range("E15").Value = "(0.00000)"
